# Fruta / hortaliza, semanal
#
# The data rows (2-39) describing each "Tuna" market transaction were
# reshuffled: for every row, the columns D (Fecha) and L..T (Calidad,
# Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg / unidad) moved to
# a different row while columns A,B,C,E..K (market/product identifiers,
# identical on every row) stayed put.
#
# Build the permutation once (new row number -> old row number the data
# came from) and copy the snapshot of the "before" values across.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

$mapping = @{
    2 = 19; 3 = 29; 4 = 35; 5 = 36; 6 = 30; 7 = 31; 8 = 21; 9 = 37; 10 = 26;
    11 = 20; 12 = 16; 13 = 3; 14 = 4; 15 = 25; 16 = 12; 17 = 32; 18 = 33;
    19 = 34; 20 = 38; 21 = 39; 22 = 22; 23 = 6; 24 = 9; 25 = 13; 26 = 14;
    27 = 15; 28 = 17; 29 = 18; 30 = 10; 31 = 11; 32 = 28; 33 = 23; 34 = 24;
    35 = 27; 36 = 5; 37 = 7; 38 = 8; 39 = 2
}

# 1) Snapshot every source cell's value before any writes happen.
#    (Value2 is used for the read because it returns the plain scalar;
#    the bare Value getter resolves to the property accessor object in
#    this host when read outside of a direct comparison/print context.)
$snapshot = @{}
foreach ($r in 2..39) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the permuted values back out.
foreach ($r in 2..39) {
    $srcRow = $mapping[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
